$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header K1 "Typ", styled like the other header cells (bold, bordered)
$ws.Range("K1").Value = "Typ"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null

# Row 2
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "95.45,17.63,32.18"
$ws.Range("K2").Value = "Inside"

# Row 3
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "37.87,13.94,13.94"
$ws.Range("K3").Value = "Inside"

# Row 4
$ws.Range("C4").Value = "351.23,49.23,19.1,19.1"
$ws.Range("K4").Value = "Outside"

# Row 5
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "193.89,15.23,44.89,15.23"
$ws.Range("K5").Value = "Inside"

# Row 6
$ws.Range("C6").Value = "60.0,30.0"
$ws.Range("D6").Value = ""
$ws.Range("K6").Value = "Outside"

# Row 7
$ws.Range("C7").Value = "60.0,30.0"
$ws.Range("D7").Value = ""
$ws.Range("K7").Value = "Outside"

# Row 8
$ws.Range("C8").Value = "70.0,49.96,17.96"
$ws.Range("D8").Value = ""
$ws.Range("K8").Value = "Outside"

# Row 9
$ws.Range("C9").Value = "200.0,50.0"
$ws.Range("D9").Value = ""
$ws.Range("K9").Value = "Outside"
